$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48:59 down to 49:60.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly data point.
$ws.Cells.Item(48, 1).Value = 1
$ws.Cells.Item(48, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value = 44704
$ws.Cells.Item(48, 5).Value = 15
$ws.Cells.Item(48, 6).Value = 100112031
$ws.Cells.Item(48, 7).Value = "Poroto verde"
$ws.Cells.Item(48, 8).Value = "Magnum"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 200
$ws.Cells.Item(48, 11).Value = 19000
$ws.Cells.Item(48, 12).Value = 20000
$ws.Cells.Item(48, 13).Value = 19500
$ws.Cells.Item(48, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value = "Perú"
$ws.Cells.Item(48, 16).Value = 780
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
